# Update Match Bank Statement
# Row 2: change date, payee, description stays "type"/"Receive", but amount/balance doubled,
#        date moves from 2020-02-03 to 2020-01-02, description shortened to "Pendapatan".
# Row 3: the old "PLN / Biaya Listrik Januari 2020 / Spend" transaction is removed (cells cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (updated transaction) ---
$ws.Range("A2").Value = 43832
# Leading apostrophe preserves the cell's existing "quote prefix" text style (s="3")
# instead of Excel resetting it to the plain bordered style when the value is reassigned.
$ws.Range("B2").Value = "'ILO"
$ws.Range("C2").Value = "Pendapatan"
$ws.Range("D2").Value = "Receive"
$ws.Range("E2").Value = 250000000
$ws.Range("F2").Value = 250000000

# --- Row 3 (transaction removed, keep empty/blank cells with their formatting) ---
$ws.Range("A3:F3").ClearContents()
